# Actualización desde MV -datos-
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing row 34 ("01-01-2021") values that changed with the data refresh
$ws.Range("H34").Value = 12484
$ws.Range("J34").Value = 12483
$ws.Range("K34").Value = 742
$ws.Range("M34").Value = 706
$ws.Range("N34").Value = 47915
$ws.Range("P34").Value = 47915
$ws.Range("T34").Value = 94308
$ws.Range("V34").Value = 85490

# Add the new row 35 ("01-04-2021") with the latest quarter's data.
# "01-04-2021" looks like a date, so a direct Range.Value assignment would get
# auto-converted to a date serial by Excel. Build it as a text formula result
# in a scratch cell, then paste-special just the value into A35 so it lands
# as plain text (matching the other "Serie" entries like "01-01-2021") without
# picking up a number-format style.
$ws.Cells.Item(36, 1).Value = "=T(""01-04-2021"")"
$ws.Range("A36").Copy()
$ws.Range("A35").PasteSpecial(-4163)
$ws.Range("A36").Clear()
$ws.Range("B35").Value = 36834
$ws.Range("C35").Value = 8807
$ws.Range("D35").Value = 28027
$ws.Range("E35").Value = 12
$ws.Range("F35").Value = 12
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 12328
$ws.Range("I35").Value = 1
$ws.Range("J35").Value = 12328
$ws.Range("K35").Value = 683
$ws.Range("L35").Value = 36
$ws.Range("M35").Value = 648
$ws.Range("N35").Value = 49620
$ws.Range("O35").Value = 2
$ws.Range("P35").Value = 49618
$ws.Range("Q35").Value = 209
$ws.Range("R35").Value = 0
$ws.Range("S35").Value = 209
$ws.Range("T35").Value = 99687
$ws.Range("U35").Value = 8857
$ws.Range("V35").Value = 90830
